$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Onto-shape with spin docs")
$ws.Activate()

# New "test category" tag column (D) added next to the existing OWL/SHACL mapping rows.
# These cells use a small monospace (Consolas) font in a dark gray, matching the
# style used elsewhere in the workbook to flag which automated test covers each row.
$tagFontName = "Consolas"
$tagFontSize = 12
$tagFontColor = 3025188  # RGB(0x24, 0x29, 0x2E) packed as BGR for OLE_COLOR

function Set-TestTag($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $text
    $rng.Font.Name = $tagFontName
    $rng.Font.Size = $tagFontSize
    $rng.Font.Color = $tagFontColor
}

# Cardinality-related rows (owl:maxCardinality/minCardinality/cardinality <-> sh:*Count)
Set-TestTag "D5" "CardinalityTest"
Set-TestTag "D6" "CardinalityTest"
Set-TestTag "D7" "CardinalityTest"

# owl:FunctionalProperty <-> sh:maxCount 1 - keeps the default cell style
$ws.Range("D8").Value = "PropertyCharacteristicsTest"

# owl:hasValue <-> sh:hasValue
Set-TestTag "D36" "FilledInformationTest"

# Annotation-property rows (rdfs:comment/label/seeAlso/isDefinedBy)
Set-TestTag "D50" "AnnotationPropertiesTest"
Set-TestTag "D51" "AnnotationPropertiesTest"
Set-TestTag "D52" "AnnotationPropertiesTest"
Set-TestTag "D53" "AnnotationPropertiesTest"

# Give the new column a sensible width now that it holds content.
$ws.Columns.Item(4).ColumnWidth = 27.1640625

# Restore the selection to the area that was being edited.
[void]$ws.Range("C13").Select()
